$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(682).Delete()
